$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3301
$ws.Range("F5").Value = 1351
$ws.Range("F7").Value = 0
$ws.Range("F9").Value = 50
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 452
$ws.Range("F12").Value = 78
$ws.Range("F13").Value = 76
$ws.Range("F14").Value = 274
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 10545
$ws.Range("F20").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F27").Value = 176
$ws.Range("F28").Value = 156
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 44
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 877
$ws.Range("F36").Value = 29
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 1227
$ws.Range("F41").Value = 165
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 63
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 274
$ws.Range("F47").Value = 99
$ws.Range("F48").Value = 82
$ws.Range("F49").Value = 0
$ws.Range("F50").Value = 0

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 186
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 0

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3301
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 386
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 8368
$ws.Range("F12").Value = 452
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 77
$ws.Range("F16").Value = 274
$ws.Range("F17").Value = 305
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 340
$ws.Range("F20").Value = 10546
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 176
$ws.Range("F29").Value = 156
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F34").Value = 877
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 2569
$ws.Range("F38").Value = 3012
$ws.Range("F39").Value = 1227
$ws.Range("F40").Value = 165
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 63
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 274
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 82
$ws.Range("F49").Value = 0
$ws.Range("F50").Value = 66
